$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.298.05'
$ws.Range("E2").Value = '  -0.08%  '
$ws.Range("D3").Value = '2.604.15'
$ws.Range("E3").Value = '  +0.47%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '''538.35'
$ws.Range("E5").Value = '  +3.48%  '
$ws.Range("D6").Value = '''140.49'
$ws.Range("E6").Value = '  +1.07%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '''0.569'
$ws.Range("E8").Value = '  +0.42%  '
$ws.Range("D9").Value = '2.609.80'
$ws.Range("E9").Value = '  -0.06%  '
$ws.Range("D10").Value = '''6.46'
$ws.Range("E10").Value = '  -0.60%  '
$ws.Range("E11").Value = '  +1.46%  '
$ws.Range("E12").Value = '  +1.17%  '
$ws.Range("D14").Value = '3.067.31'
$ws.Range("E14").Value = '  +0.32%  '
$ws.Range("D15").Value = '59.219.48'
$ws.Range("E15").Value = '  -0.15%  '
$ws.Range("D16").Value = '''20.53'
$ws.Range("E16").Value = '  +0.78%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '2.616.69'
$ws.Range("E17").Value = '  -0.10%  '
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").Value = '''0.0000134'
$ws.Range("E18").Value = '  +0.73%  '
$ws.Range("D19").Value = '''342.82'
$ws.Range("E19").Value = '  +1.12%  '
$ws.Range("E20").Value = '  +0.49%  '
$ws.Range("D21").Value = '''10.09'
$ws.Range("E21").Value = '  -0.77%  '
$ws.Range("D22").Value = '''6.39'
$ws.Range("E22").Value = '  -1.61%  '
$ws.Range("D23").Value = '''1.00'
$ws.Range("E23").Value = '  +0.21%  '
$ws.Range("D24").Value = '''67.56'
$ws.Range("E24").Value = '  +1.80%  '
$ws.Range("E25").Value = '  -0.56%  '
$ws.Range("E26").Value = '  +1.11%  '
$ws.Range("E27").Value = '  +0.26%  '
$ws.Range("D28").Value = '''7.20'
$ws.Range("E28").Value = '  +2.47%  '
$ws.Range("E29").Value = '  +0.11%  '
$ws.Range("D30").Value = '0.0₃0738'
$ws.Range("E30").Value = '  +2.02%  '
$ws.Range("E31").Value = '  +5.39%  '
$ws.Range("D32").Value = '''5.83'
$ws.Range("E32").Value = '  -2.15%  '
$ws.Range("D33").Value = '''18.79'
$ws.Range("E33").Value = '  -0.15%  '
$ws.Range("D34").Value = '''149.36'
$ws.Range("E34").Value = '  +0.03%  '
$ws.Range("D35").Value = '''3.97'
$ws.Range("E35").Value = '  -0.49%  '
$ws.Range("E36").Value = '  -1.03%  '
$ws.Range("D37").Value = '''36.92'
$ws.Range("E37").Value = '  +1.71%  '
$ws.Range("D38").Value = '''1.47'
$ws.Range("E38").Value = '  +1.04%  '
$ws.Range("D39").Value = '''0.836'
$ws.Range("E39").Value = '  +1.16%  '
$ws.Range("D40").Value = '''0.827'
$ws.Range("E40").Value = '  +0.53%  '
$ws.Range("D41").Value = '''3.55'
$ws.Range("E41").Value = '  +0.24%  '
$ws.Range("E42").Value = '  +0.12%  '
$ws.Range("D43").Value = '''274.37'
$ws.Range("E43").Value = '  -0.08%  '
$ws.Range("D44").Value = '''0.596'
$ws.Range("E44").Value = '  +0.81%  '
$ws.Range("D45").Value = '''10.75'
$ws.Range("E45").Value = '  +0.04%  '
$ws.Range("E46").Value = '  +1.08%  '
$ws.Range("D47").Value = '''0.0523'
$ws.Range("E47").Value = '  +0.50%  '
$ws.Range("D48").Value = '1.946.36'
$ws.Range("E48").Value = '  -1.87%  '
$ws.Range("E49").Value = '  +1.17%  '
$ws.Range("D50").Value = '''18.28'
$ws.Range("E50").Value = '  +1.37%  '
$ws.Range("D51").Value = '''4.50'
$ws.Range("E51").Value = '  +0.48%  '
